$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear C2 and B3 (they no longer hold a value)
$ws.Range("C2").ClearContents()
$ws.Range("B3").ClearContents()

# Remove rows 4-7 entirely (data + used range shrinks to A1:C3)
$ws.Range("A4:C7").ClearContents()

# Update the selected cell to match the saved view state
$ws.Range("A6").Select()
